# Adds a new week of data (week of 2022-03-17, serial date 44637) for
# "Femacal de La Calera - Membrillo" by inserting two new rows right
# before the existing row 58 block, shifting all rows from 58 downward
# down by two (old row 58 -> new row 60, ..., old row 75 -> new row 77).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 58, pushing the rest of the
# table (old rows 58..75) down to rows 60..77.
$ws.Range("A58:A59").EntireRow.Insert()

# --- New row 58: Calidad "Especial" ---
$ws.Range("A58").Value = 3
$ws.Range("B58").Value = "Femacal de La Calera"
$ws.Range("C58").Value = "Coquimbo"
$ws.Range("D58").Value = 44637
$ws.Range("E58").Value = 5
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100104
$ws.Range("H58").Value = "Frutos de pepita"
$ws.Range("I58").Value = 100104003
$ws.Range("J58").Value = "Membrillo"
$ws.Range("K58").Value = "Champion"
$ws.Range("L58").Value = "Especial"
$ws.Range("M58").Value = 75
$ws.Range("N58").Value = 18000
$ws.Range("O58").Value = 18000
$ws.Range("P58").Value = 18000
$ws.Range("Q58").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R58").Value = "Región de O'Higgins"
$ws.Range("S58").Value = 1000
$ws.Range("T58").Value = 18

# --- New row 59: Calidad "Primera" ---
$ws.Range("A59").Value = 3
$ws.Range("B59").Value = "Femacal de La Calera"
$ws.Range("C59").Value = "Coquimbo"
$ws.Range("D59").Value = 44637
$ws.Range("E59").Value = 5
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100104
$ws.Range("H59").Value = "Frutos de pepita"
$ws.Range("I59").Value = 100104003
$ws.Range("J59").Value = "Membrillo"
$ws.Range("K59").Value = "Champion"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 70
$ws.Range("N59").Value = 16000
$ws.Range("O59").Value = 16000
$ws.Range("P59").Value = 16000
$ws.Range("Q59").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R59").Value = "Región de O'Higgins"
$ws.Range("S59").Value = 889
$ws.Range("T59").Value = 18
